# Apply the "data_moments" table update:
#  - Add two new header columns (D: FIRE+AR, E: FIRE+SV)
#  - Populate the new columns with text/numeric "moment" markers per row
#  - Row labels (column A) and B/C numeric values are otherwise unchanged,
#    except B2/C2 which move from 2.146/2.251 to (numerically) 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells, copying the style of the existing headers (B1/C1) ---
$ws.Range("D1").Value2 = 'FIRE+AR'
$ws.Range("E1").Value2 = 'FIRE+SV'

$ws.Range("B1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Updated numeric values in B2 / C2 ---
$ws.Range("B2").Value2 = 0
$ws.Range("C2").Value2 = 0

# --- New column D values (row 2 .. row 13) ---
$ws.Range("D2").Value2 = 0
$ws.Range("D3").Value2 = '$\sigma^2/(1-\rho^2)$'
$ws.Range("D4").Value2 = '$\rho\sigma^2/(1-\rho^2)$'
$ws.Range("D5").Value2 = 0
$ws.Range("D6").Value2 = '$\sigma^2$'
$ws.Range("D7").Value2 = 0
$ws.Range("D8").Value2 = 0
$ws.Range("D9").Value2 = 0
$ws.Range("D10").Value2 = 0
$ws.Range("D11").Value2 = '$\sigma^2$'
$ws.Range("D12").Value2 = 0
$ws.Range("D13").Value2 = 0

# --- New column E values (row 2 .. row 13) ---
$ws.Range("E2").Value2 = 'N/A'
$ws.Range("E3").Value2 = 'N/A'
$ws.Range("E4").Value2 = 'N/A'
$ws.Range("E5").Value2 = 0
$ws.Range("E6").Value2 = '$\bar\sigma^2_{\eta}+\bar\sigma^2_{\epsilon}$'
$ws.Range("E7").Value2 = 0
$ws.Range("E8").Value2 = 0
$ws.Range("E9").Value2 = 0
$ws.Range("E10").Value2 = 0
$ws.Range("E11").Value2 = '$\bar\sigma^2_{\eta}+\bar\sigma^2_{\epsilon}$'
$ws.Range("E12").Value2 = '>0'
$ws.Range("E13").Value2 = '>0'
